$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of filenames appended below the existing list (A16:A21).
# Entered bottom-up so the shared-string table records them in
# alphabetical order (matching the saved file's sharedStrings.xml).
$ws.Range("A21").Value = "SV 1127 Flatbush Avenue.xlsx"
$ws.Range("A20").Value = "SV 148 West Street.xlsx"
$ws.Range("A19").Value = "SV 2501 Pitkin Avenue.xlsx"
$ws.Range("A18").Value = "SV 2702 West 15th Street.xlsx"
$ws.Range("A17").Value = "SV 432-436 Keap Street.xlsx"
$ws.Range("A16").Value = "SV 929 Atlantic Avenue.xlsx"

# Update the selected cell to match the saved view state
$ws.Range("D18").Select()
